$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ---------------------------------------------------------------------------
# 1) Merge paragraph 1 ("...clinical researchers.") with paragraph 2
#    ("They asked ... members.") joined by a single space, as three runs:
#      run1 = original paragraph-1 text (untouched)
#      run2 = " "
#      run3 = original paragraph-2 text (untouched)
# ---------------------------------------------------------------------------

# 1a) Insert a new run containing just a space right before the paragraph
#     mark that currently separates paragraph 1 and paragraph 2. Replacing
#     the (1-char) paragraph-mark selection with a <w:p> wrapper keeps the
#     paragraph break itself intact but adds the new run ahead of it,
#     without disturbing the existing run's text/formatting.
$rng = $d.Content
$rng.Find.Execute("clinical researchers.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.MoveEnd(1, 1)
$xml = '<w:p xmlns:w="' + $wNs + '"><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$rng.InsertXML($xml)

# 1b) Now delete the (new) paragraph mark that still separates the two
#     paragraphs, which merges them into one paragraph without touching any
#     run text.
$rng2 = $d.Content
$rng2.Find.Execute("clinical researchers. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Collapse(0)
$rng2.MoveEnd(1, 1)
$rng2.Delete()

# ---------------------------------------------------------------------------
# 2) Insert the two brand-new paragraphs (second & third data-science
#    groups) right after the merged paragraph 1.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertAfter("The second group is Algorithmic and machine learning group. They basically focusing on understanding mathematic for improving machining learning models. I learned what mathematic do I need for long standing data science.")
$r2.InsertParagraphAfter()

$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertAfter("The third group is Kaggle Austin. They introduce 2 persons every meeting to present their projects. Also, they encourage us to use Kaggle and upload our data sets. For this group, I could not understand their goal. It seems that they advertise to Kaggle website.")
$r3.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 3) Insert the fourth (new) paragraph about Ascension information, split
#    across two runs exactly like the source diff, and carry the _GoBack
#    bookmark down to the very end of the document (it always trails the
#    last edit).
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$r4.Collapse(0)
$r4.InsertAfter("One of the groups I met and learned from them is Ascension information. I met them in the second interview for data science job position. They want the data scientist to have some experience and will to work, beside data scientist, as data engineer. They skills they look for is big data wrangling using Hadoop and spark.")

$rng4 = $d.Content
$rng4.Find.Execute("Hadoop and spark.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$xml4 = '<w:p xmlns:w="' + $wNs + '"><w:r><w:t>Hadoop and spark.</w:t></w:r><w:r><w:t xml:space="preserve"> I learned the principles of these two technologies and will continue with them after finishing my course.</w:t></w:r></w:p>'
$rng4.InsertXML($xml4)

# ---------------------------------------------------------------------------
# 4) Move the (hidden) _GoBack bookmark to the very end of the document, so
#    it once again marks the spot right after the final edit.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $endRng)

Write-Host "Done. Paragraphs:" $d.Paragraphs.Count
